$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to insert right after the header row (row 1),
# pushing the existing data down by 4 rows.
$newTopRows = @(
    @(-0.1621847152709961, -0.07956519722938531, -0.5577199459075928),
    @(-0.1879937797784805, -0.1305724531412124,  -0.08170322328805921),
    @(-0.0765108689665794, -0.2229658216238021,  -0.0765108689665794),
    @(-0.0930042341351509, 0.1044579595327377,   -0.295353353023529)
)

# Insert 4 new rows before row 2 (shifts existing data rows 2-21 down to 6-25)
$insertRange = $ws.Range("A2:C5")
$insertRange.EntireRow.Insert() | Out-Null

# The inserted rows pick up formatting from the row above (the bold header);
# clear that so the new rows look like ordinary (unstyled) data rows.
$ws.Range("A2:C5").ClearFormats() | Out-Null

# Fill the newly inserted rows with their values
for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTopRows[$i][2]
}

# New rows appended at the bottom of the data (rows 26-31)
$newBottomRows = @(
    @(-0.0247400421649217, 0.0786489024758338, -0.0580321997404098),
    @(-0.026419922709465,  0.0665843114256858, -0.07590000331401819),
    @(0.0235183127224445,  0.012980886735022,  -0.0116064399480819),
    @(-0.009315694682300001, 0.0103847095742821, -0.0035124751739203),
    @(-0.0279470849782228, 0.0250454749912023, -0.00335975876078),
    @(-0.0271835029125213, 0.0594066455960273, -0.0343611687421798)
)

for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = 26 + $i
    $ws.Cells.Item($r, 1).Value = $newBottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottomRows[$i][2]
}
